# Updated cryptos list (price/volume refresh + two row re-rankings).
# Price/Volume(1h) cells are plain text (e.g. "1.000", dotted thousands
# like "29.385.40"), so any cell whose new value Excel would otherwise
# auto-coerce to a number is first forced to text format ("@") before
# the assignment, preserving exact formatting (trailing zeros, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.385.40"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.849.81"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.26"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6286"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07627"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2905"
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.75"
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.035"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6786"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.29"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.167"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "29.407.97"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.00"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.36"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.500"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "158.68"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1388"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.406"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.71"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.387"
$ws.Range("E27").Value = "  +6.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.463"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05603"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.110"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.070"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.164"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.838"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7009"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.580"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01807"
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").Value = "1.232.36"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.385"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9026"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.54"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.11"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.218"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000119"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4013"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.997"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.681"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1134"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05701"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4630"
$ws.Range("E51").Value = "  +0.09%  "
